$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).ColumnWidth = 24.54296875
$ws.Columns.Item(2).ColumnWidth = 29.453125
